$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data (2022-07-11) is inserted above the most recent
# historical rows, pushing the existing rows 47-49 down to 48-50.
$ws.Rows("47:47").Insert()

$ws.Range("A47").Value = 8
$ws.Range("B47").Value = "Terminal La Palmera de La Serena"
$ws.Range("C47").Value = "Coquimbo"
$ws.Range("D47").Value = 44753
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = 100114007
$ws.Range("G47").Value = "Jengibre"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 400
$ws.Range("K47").Value = 14500
$ws.Range("L47").Value = 15000
$ws.Range("M47").Value = 14750
$ws.Range("N47").Value = "`$/caja 13 kilos"
$ws.Range("O47").Value = "Perú"
$ws.Range("P47").Value = 1135
$ws.Range("Q47").Value = 13
$ws.Range("R47").Value = "Hortaliza"
